# Generate Report for Handback
# Adds a new handback record for c4ec24ab-e82a-46b2-9692-3a3ff5ecfe3a
# (in sync with en-US) as row 4 on the "Overview", "zh-cn" and "de-de"
# worksheets, mirroring the existing 60f85624-a633-455a-b979-4449fc044161
# ("Include" / in-sync) row.

$wb = $excel.ActiveWorkbook

$uuid = "c4ec24ab-e82a-46b2-9692-3a3ff5ecfe3a"
$xlfHash = "4a140e5060350133549b5a58d04cbcfc2f3d3631"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$uuid.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b2f2b14cfb06b4f6e5c8d6c2f7f0b2d9c4a6e8f1/e2e/$uuid.md",
    "",
    "",
    "$uuid.md"
) | Out-Null

$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------
# Locale sheets: "zh-cn" and "de-de"
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; HandoffTime = "2016-03-11 08:31:11"; HandbackTime = "2016-03-11 08:31:56";
       HandoffSha = "c1e6f0a7b6d4e2c8f9a3b5d7e1c4f6a8b0d2e4f6"; MdSha2 = "9f3a5c7e1b4d6f8a0c2e4f6a8b0d2e4f6a8b0d2e";
       XlfHandoffSha = "7a4c6e8f0b2d4f6a8c0e2f4a6b8d0e2f4a6b8d0e"; XlfHandbackSha = "3e5a7c9f1b3d5f7a9c1e3f5a7c9e1f3a5c7e9f1b" },
    @{ Sheet = "de-de"; HandoffTime = "2016-03-11 08:31:22"; HandbackTime = "2016-03-11 08:32:12";
       HandoffSha = "d2f7a1b8c6e4d2f0a8c6e4d2f0a8c6e4d2f0a8c6"; MdSha2 = "1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b";
       XlfHandoffSha = "6b8d0f2a4c6e8b0d2f4a6c8e0b2d4f6a8c0e2f4a"; XlfHandbackSha = "4f6a8c0e2b4d6f8a0c2e4f6a8c0e2f4a6c8e0b2d" }
)

foreach ($loc in $locales) {
    $sheetName = $loc.Sheet
    $ws = $wb.Worksheets.Item($sheetName)

    $mdDisplay = "$uuid.md"
    $xlfName = "$uuid.$xlfHash.$sheetName.xlf"

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$($loc.MdSha2)/e2e/$mdDisplay"
    $handoffXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($loc.XlfHandoffSha)/ol-handoff/OpenLocalizationTestOrg/oltest.$sheetName/yuwzho/ht/$xlfName"
    $targetMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest.$sheetName/blob/$($loc.HandoffSha)/e2e/$mdDisplay"
    $handbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($loc.XlfHandbackSha)/ol-handback/OpenLocalizationTestOrg/oltest.$sheetName/yuwzho/ht/$xlfName"

    # A4: Source File Name (md)
    $ws.Range("A4").Value = $mdDisplay
    $ws.Hyperlinks.Add($ws.Range("A4"), $mdUrl, "", "", $mdDisplay) | Out-Null

    # B4: File Extension (".md") - links to the same source file
    $ws.Range("B4").Value = ".md"
    $ws.Hyperlinks.Add($ws.Range("B4"), $mdUrl, "", "", ".md") | Out-Null

    # C4: Status
    $ws.Range("C4").Value = $statusInSync

    # D4: Correspond Handoff File (xlf)
    $ws.Range("D4").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("D4"), $handoffXlfUrl, "", "", $xlfName) | Out-Null

    # E4: Correspond Handoff Datetime
    $ws.Range("E4").Value = $loc.HandoffTime

    # F4: Target File (md)
    $ws.Range("F4").Value = $mdDisplay
    $ws.Hyperlinks.Add($ws.Range("F4"), $targetMdUrl, "", "", $mdDisplay) | Out-Null

    # G4: Correspond Handback File (xlf)
    $ws.Range("G4").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("G4"), $handbackXlfUrl, "", "", $xlfName) | Out-Null

    # H4: Correspond Handback DateTime
    $ws.Range("H4").Value = $loc.HandbackTime

    # I4: Handoff Reason
    $ws.Range("I4").Value = "Include"
}
